# Update view-count figures (column F) on the "展览" and "全部类型" sheets
# to match freshly generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 88
$wsExhibit.Range("F10").Value = 15790
$wsExhibit.Range("F14").Value = 6197
$wsExhibit.Range("F19").Value = 120
$wsExhibit.Range("F27").Value = 868
$wsExhibit.Range("F31").Value = 11092
$wsExhibit.Range("F34").Value = 127
$wsExhibit.Range("F35").Value = 177

# Sheet "全部类型" (4th sheet)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 88
$wsAll.Range("F11").Value = 15790
$wsAll.Range("F15").Value = 6197
$wsAll.Range("F20").Value = 120
$wsAll.Range("F28").Value = 868
$wsAll.Range("F33").Value = 11092
$wsAll.Range("F36").Value = 127
$wsAll.Range("F37").Value = 177
